$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "is_profile_data_fetched" column (B) for the rows that were
# previously missing it (rows 25-132), then append one more scraped
# username (row 133) with both flags marked "done".
$ws.Range("B25:B132").Value = "done"
$ws.Range("A133").Value = "dr.rakshita_singh"
$ws.Range("B133").Value = "done"

# Widen the columns now that header / data text is wider.
$ws.Columns.Item(1).ColumnWidth = 21.7
$ws.Columns.Item(2).ColumnWidth = 22.0
$ws.Columns.Item(3).ColumnWidth = 23.0

# Turn the populated range into a proper Excel Table.
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:C133"), $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = "TableStyleMedium9"

# Restore the view's selection to where the user left off.
$ws.Range("M122").Select()

Write-Output "edit applied"
